$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.338.19"
$ws.Range("E2").Value = "  +1.46%  "

$ws.Range("D3").Value = "1.886.62"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Formula = "'245.64"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("E6").Value = "  +1.05%  "

$ws.Range("D8").Formula = "'42.86"
$ws.Range("E8").Value = "  +2.04%  "

$ws.Range("E9").Value = "  +2.26%  "

$ws.Range("D10").Formula = "'55.18"
$ws.Range("E10").Value = "  +8.11%  "

$ws.Range("D11").Formula = "'0.0743"
$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("E12").Value = "  +1.35%  "

$ws.Range("D13").Formula = "'13.82"
$ws.Range("E13").Value = "  +7.18%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Formula = "'0.776"
$ws.Range("E14").Value = "  +8.71%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.159.93"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("D16").Formula = "'4.99"
$ws.Range("E16").Value = "  +2.53%  "

$ws.Range("D17").Value = "1.897.00"
$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("D18").Value = "35.287.43"
$ws.Range("E18").Value = "  +1.32%  "

$ws.Range("D19").Formula = "'73.35"
$ws.Range("E19").Value = "  +0.98%  "

$ws.Range("D20").Value = "0.0₃0824"
$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("D21").Formula = "'244.67"
$ws.Range("E21").Value = "  +0.67%  "

$ws.Range("D22").Formula = "'12.81"
$ws.Range("E22").Value = "  +1.32%  "

$ws.Range("E23").Value = "  +4.82%  "

$ws.Range("E24").Value = "  +7.79%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").Formula = "'2.16"
$ws.Range("E26").Value = "  -2.57%  "

$ws.Range("D27").Formula = "'167.29"
$ws.Range("E27").Value = "  +1.41%  "

$ws.Range("D28").Formula = "'8.54"
$ws.Range("E28").Value = "  +2.07%  "

$ws.Range("D29").Formula = "'18.27"
$ws.Range("E29").Value = "  +0.51%  "

$ws.Range("E30").Value = "  +0.67%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Formula = "'4.28"
$ws.Range("E31").Value = "  +1.88%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Formula = "'0.0594"
$ws.Range("E32").Value = "  +3.00%  "

$ws.Range("E33").Value = "  +24.36%  "

$ws.Range("D34").Formula = "'4.19"
$ws.Range("E34").Value = "  +1.20%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("E36").Value = "  -13.86%  "

$ws.Range("D37").Formula = "'0.852"
$ws.Range("E37").Value = "  +3.28%  "

$ws.Range("D38").Formula = "'1.94"
$ws.Range("E38").Value = "  -2.17%  "

$ws.Range("D39").Formula = "'0.0720"
$ws.Range("E39").Value = "  +8.66%  "

$ws.Range("D40").Formula = "'0.0222"
$ws.Range("E40").Value = "  +6.09%  "

$ws.Range("D41").Formula = "'98.13"
$ws.Range("E41").Value = "  +0.82%  "

$ws.Range("D42").Formula = "'17.16"
$ws.Range("E42").Value = "  +0.35%  "

$ws.Range("E43").Value = "  -1.15%  "

$ws.Range("D44").Formula = "'13.67"
$ws.Range("E44").Value = "  +15.33%  "

$ws.Range("D45").Value = "1.326.90"
$ws.Range("E45").Value = "  +3.34%  "

$ws.Range("E46").Value = "  +2.89%  "

$ws.Range("D47").Formula = "'0.0809"
$ws.Range("E47").Value = "  +0.49%  "

$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("D50").Formula = "'6.28"
$ws.Range("E50").Value = "  -2.07%  "

$ws.Range("D51").Value = "2.059.36"
$ws.Range("E51").Value = "  +0.13%  "

